$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "I2" = 0.977669497583861
    "J2" = 0.977669497583861
    "M2" = 1.090710333333333
    "N2" = 3.272131
    "O2" = 0.0488470045579656
    "P2" = 0.0488470045579656
    "Q2" = 8.196464559381667
    "R2" = 73.768181034435
    "S2" = 0.0477562264046628
    "T2" = 0.0477562264046628
    "I3" = 0.977669497583861
    "J3" = 0.977669497583861
    "O3" = 0.7616320856558244
    "P3" = 0.7616320856558244
    "S3" = 0.7446244585268781
    "T3" = 0.7446244585268781
    "I4" = 0.977669497583861
    "J4" = 0.977669497583861
    "M4" = 0.740281
    "N4" = 2.220843
    "O4" = 0.03315317392351528
    "P4" = 0.03315317392351528
    "Q4" = 5.563059957395
    "R4" = 50.067539616555
    "S4" = 0.03241284689311354
    "T4" = 0.03241284689311354
    "I5" = 0.977669497583861
    "J5" = 0.977669497583861
    "M5" = 2.784013333333333
    "N5" = 8.352039999999999
    "O5" = 0.1246808688124989
    "P5" = 0.1246808688124989
    "Q5" = 20.92128947726666
    "R5" = 188.2916052954
    "S5" = 0.1218966823702351
    "T5" = 0.1218966823702351
    "I6" = 0.977669497583861
    "J6" = 0.977669497583861
    "M6" = 0.2710316666666667
    "N6" = 0.813095
    "O6" = 0.01213803945228936
    "P6" = 0.01213803945228936
    "Q6" = 2.036747413508333
    "R6" = 18.330726721575
    "S6" = 0.01186699093297282
    "T6" = 0.01186699093297282
    "I7" = 0.977669497583861
    "J7" = 0.977669497583861
    "M7" = 0.436508
    "N7" = 1.309524
    "O7" = 0.01954882759790648
    "P7" = 0.01954882759790648
    "Q7" = 3.280268135859999
    "R7" = 29.52241322274
    "S7" = 0.01911229245599875
    "T7" = 0.01911229245599875
    "E8" = 2
    "F8" = 0.6666666666666666
    "G8" = 0.171642
    "H8" = 0.514926
    "I8" = 0.02233050241613897
    "J8" = 0.02233050241613898
    "M8" = 1.090710333333333
    "N8" = 3.272131
    "O8" = 0.0488470045579656
    "P8" = 0.0488470045579656
    "Q8" = 0.187211703034
    "R8" = 1.684905327306
    "S8" = 0.001090778153302802
    "T8" = 0.001090778153302802
    "E9" = 2
    "F9" = 0.6666666666666666
    "G9" = 0.171642
    "H9" = 0.514926
    "I9" = 0.02233050241613897
    "J9" = 0.02233050241613898
    "O9" = 0.7616320856558244
    "P9" = 0.7616320856558244
    "Q9" = 2.919041630726
    "R9" = 26.271374676534
    "S9" = 0.01700762712894635
    "T9" = 0.01700762712894635
    "E10" = 2
    "F10" = 0.6666666666666666
    "G10" = 0.171642
    "H10" = 0.514926
    "I10" = 0.02233050241613897
    "J10" = 0.02233050241613898
    "M10" = 0.740281
    "N10" = 2.220843
    "O10" = 0.03315317392351528
    "P10" = 0.03315317392351528
    "Q10" = 0.127063311402
    "R10" = 1.143569802618
    "S10" = 0.0007403270304017335
    "T10" = 0.0007403270304017337
    "E11" = 2
    "F11" = 0.6666666666666666
    "G11" = 0.171642
    "H11" = 0.514926
    "I11" = 0.02233050241613897
    "J11" = 0.02233050241613898
    "M11" = 2.784013333333333
    "N11" = 8.352039999999999
    "O11" = 0.1246808688124989
    "P11" = 0.1246808688124989
    "Q11" = 0.4778536165599999
    "R11" = 4.300682549039999
    "S11" = 0.002784186442263813
    "T11" = 0.002784186442263814
    "E12" = 2
    "F12" = 0.6666666666666666
    "G12" = 0.171642
    "H12" = 0.514926
    "I12" = 0.02233050241613897
    "J12" = 0.02233050241613898
    "M12" = 0.2710316666666667
    "N12" = 0.813095
    "O12" = 0.01213803945228936
    "P12" = 0.01213803945228936
    "Q12" = 0.04652041733
    "R12" = 0.41868375597
    "S12" = 0.0002710485193165377
    "T12" = 0.0002710485193165378
    "E13" = 2
    "F13" = 0.6666666666666666
    "G13" = 0.171642
    "H13" = 0.514926
    "I13" = 0.02233050241613897
    "J13" = 0.02233050241613898
    "M13" = 0.436508
    "N13" = 1.309524
    "O13" = 0.01954882759790648
    "P13" = 0.01954882759790648
    "Q13" = 0.07492310613599999
    "R13" = 0.674307955224
    "S13" = 0.0004365351419077349
    "T13" = 0.000436535141907735
}

foreach ($key in $data.Keys) {
    $ws.Range($key).Value = $data[$key]
}